# Commit: "Fruta / hortaliza, semanal"
# A new weekly price record is inserted as a new row 641 in the data table,
# pushing the previously existing rows 641-686 down to 642-687.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 641 (existing row 641 and below shift down)
$ws.Rows("641:641").Insert()

# Populate the new row with the new observation
$ws.Cells.Item(641, 1).Value = 4
$ws.Cells.Item(641, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(641, 3).Value = "Los Lagos"
$ws.Cells.Item(641, 4).Value = [DateTime]"2023-12-05"
$ws.Cells.Item(641, 5).Value = 10
$ws.Cells.Item(641, 6).Value = 100114013
$ws.Cells.Item(641, 7).Value = "Zanahoria"
$ws.Cells.Item(641, 8).Value = "Sin especificar"
$ws.Cells.Item(641, 9).Value = "Primera"
$ws.Cells.Item(641, 10).Value = 800
$ws.Cells.Item(641, 11).Value = 9000
$ws.Cells.Item(641, 12).Value = 10000
$ws.Cells.Item(641, 13).Value = 9500
$ws.Cells.Item(641, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(641, 15).Value = "Región Metropolitana"
$ws.Cells.Item(641, 16).Value = 475
$ws.Cells.Item(641, 17).Value = 20
$ws.Cells.Item(641, 18).Value = "Hortaliza"
